$d = $word.ActiveDocument

function Replace-ParagraphText {
    param(
        [string]$OldText,
        [string]$RunInnerXml
    )

    foreach ($p in $d.Paragraphs) {
        $pText = $p.Range.Text
        if ($pText.TrimEnd("`r", "`n") -eq $OldText) {
            $start = $p.Range.Start
            $end = $p.Range.End
            # Exclude the trailing paragraph mark from the range so only the
            # runs (including any leading empty run) are replaced in place;
            # the paragraph's own pPr / paragraph mark are left untouched.
            $rng = $d.Range($start, $end - 1)

            $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body><w:p>' + $RunInnerXml + '</w:p></w:body>' +
                '</w:document>' +
                '</pkg:xmlData></pkg:part></pkg:package>'

            $rng.InsertXML($xml)
            return $true
        }
    }
    return $false
}

Replace-ParagraphText `
    "Iconic historical figures portrayed as caricatures" `
    '<w:r><w:t>Interesting theme based on the London subway</w:t></w:r>'

Replace-ParagraphText `
    "Available on both PC and portable devices" `
    '<w:r><w:t>Compatible with PC and portable devices</w:t></w:r>'

Replace-ParagraphText `
    "High volatility may not be suitable for all players" `
    '<w:r><w:t>High volatility may not appeal to all players</w:t></w:r>'

Replace-ParagraphText `
    "Limited betting options compared to some other slot games" `
    '<w:r><w:t>RTP value of 95.58% is slightly below average</w:t></w:r>'

Replace-ParagraphText `
    "Read our review of Down the Rails, the London subway-themed slot game from Pragmatic Play. Play Down the Rails for free and enjoy bonus games and random features." `
    '<w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Down the Rails, a beautifully designed slot game by Pragmatic Play. Play for free and experience the London subway-themed gameplay.</w:t></w:r>'
